$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

# Fix E40/E41 to be true numbers instead of text-inline-strings
$ws.Range("E40").Value = 20
$ws.Range("E41").Value = 531344

# New breakout rows appended at the bottom (rows 42-44)
$rows = @(
    @("24/06/2024 09:44:41", 1, "BSE",    "BSE (Bombay stock exchange)",            20,     -2.55, 2494.95, 861725),
    @("24/06/2024 09:44:41", 2, "CONCOR", "Container Corporation Of India Limited", 531344, -4.02, 1047.05, 4804800),
    @("24/06/2024 09:44:41", 3, "ZEEL",   "Zee Entertainment Enterprises Limited",  505537, -2.02, 151.13,  10078997)
)

$r = 42
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    # bsecode stays a text value on new rows (matches the un-touched E40/E41 look before their fix)
    $ecell = $ws.Cells.Item($r, 5)
    $ecell.Value = "'" + $row[4]
    $ecell.Style = "Normal"
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r++
}
